# Update the "取得日時" (retrieved datetime) timestamps in column A (rows 2-6)
# on the "ランサーズ" sheet to reflect the new append time: 2025-12-30 18:37:04

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-30 18:37:04"

for ($row = 2; $row -le 6; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
